$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows containing "health" / "com.lifeomic.lifeextend" are being removed.
# In the original sheet these are rows 4, 11, 18 and 22.
# Delete them from bottom to top so earlier row numbers stay valid while
# later rows shift up, just like a user would do via Ctrl+Click row
# headers + Delete.
$ws.Range("A22:B22").EntireRow.Delete() | Out-Null
$ws.Range("A18:B18").EntireRow.Delete() | Out-Null
$ws.Range("A11:B11").EntireRow.Delete() | Out-Null
$ws.Range("A4:B4").EntireRow.Delete() | Out-Null

# Reproduce the resulting selection/active cell left behind by Excel after
# such a multi-row delete operation (originally a multi-area selection of
# rows 4, 11, 18 and 22 with A4 as the active cell).
$excel.Union($ws.Range("A4:XFD4"), $ws.Range("A11:XFD11"), $ws.Range("A18:XFD18"), $ws.Range("A22:XFD22")).Select() | Out-Null
